$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B7").Value = 100204000000.0
$ws.Range("G7").Value = 73646000000.0
$ws.Range("B8").Value = 25757000000.0
$ws.Range("B11").Value = 3704000000.0
$ws.Range("B15").Value = 3819000000.0
$ws.Range("G15").Value = 5561000000.0
$ws.Range("B17").Value = 2978000000.0
$ws.Range("B20").Value = 45781000000.0
$ws.Range("B22").Value = 25269000000.0
$ws.Range("G22").Value = 4554000000.0
$ws.Range("B24").Value = 4377000000.0
$ws.Range("G24").Value = 1701000000.0
$ws.Range("B29").Value = 1000000.0
$ws.Range("G29").Value = 688000.0
$ws.Range("B36").Value = -106857000000.0
$ws.Range("B37").Value = 28247000000.0
